$p = $ppt.ActivePresentation

# --- Slide 1: add "Đặng Minh " before "Trí" in the team member list ---
$s1 = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item(2)
$tr1 = $sh1.TextFrame.TextRange
$lastPara = $tr1.Paragraphs(5)
$nameRun = $lastPara.Runs(1)
$nameRun.InsertBefore("Đặng Minh ")

# --- Slide 3: merge "thị các hình ảnh vui, clip vui về thú " + "cưng" into one run ---
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(2)
$tr3 = $sh3.TextFrame.TextRange
$funcPara = $tr3.Paragraphs(4)
$run2 = $funcPara.Runs(2)
$run3 = $funcPara.Runs(3)
$merged = $run2.Text + $run3.Text
$run2.Text = $merged
$run3.Text = ""
